# Generate Report for Handoff
# A new handoff was generated for the "3d67bfa0-762f-44da-91ef-c03df45add61" file
# (row 4 in every sheet). This bumps that row's "Latest Handoff Date(time)" cells
# to a fresh timestamp, in the Overview summary sheet and in both per-locale
# detail sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D4").Value = "2016-03-24 09:59:59"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-24 09:59:50"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-24 09:59:59"
